$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics after trade #68 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.17   # Current Capital
$summary.Range("B4").Value = 0.16      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 68        # Total Trades
$summary.Range("B8").Value = 36        # Losing Trades
$summary.Range("B9").Value = 30.88     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: refresh MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.17     # Capital
$status.Range("D4").Value = 68         # Trades
$status.Range("E4").Value = 0.16       # P&L $
$status.Range("F4").Value = 0.17       # P&L %
$status.Range("G4").Value = 30.88      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade #68 as row 69 on both the "All Trades" and
# "MarketMaking" sheets (each keeps a full copy of every trade row).
# ---------------------------------------------------------------------------
$newRow = @{
    A = 68
    B = "2026-02-17"
    C = "15:47:58"
    D = "MarketMaking"
    E = "UP"
    F = 0.76
    G = 0.63
    H = "CLOSED"
    I = -17.1053
    J = -0.13
    K = 100.17
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 69
    $ws.Cells.Item($r, 1).Value = $newRow.A
    # Leading apostrophe forces text entry so the date-shaped string isn't
    # auto-converted into a date serial, matching the other rows above it.
    $ws.Cells.Item($r, 2).Value = "'" + $newRow.B
    $ws.Cells.Item($r, 3).Value = $newRow.C
    $ws.Cells.Item($r, 4).Value = $newRow.D
    $ws.Cells.Item($r, 5).Value = $newRow.E
    $ws.Cells.Item($r, 6).Value = $newRow.F
    $ws.Cells.Item($r, 7).Value = $newRow.G
    $ws.Cells.Item($r, 8).Value = $newRow.H
    $ws.Cells.Item($r, 9).Value = $newRow.I
    $ws.Cells.Item($r, 10).Value = $newRow.J
    $ws.Cells.Item($r, 11).Value = $newRow.K
    $ws.Cells.Item($r, 12).Value = $newRow.L
    $ws.Cells.Item($r, 13).Value = $newRow.M
    $ws.Cells.Item($r, 14).Value = $newRow.N
    $ws.Cells.Item($r, 15).Value = $newRow.O
    $ws.Cells.Item($r, 16).Value = $newRow.P
    $ws.Cells.Item($r, 17).Value = $newRow.Q
}
